# Commit: feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" (fund-holding detail, same shape as
#    "2021-Q4") right after the existing "2021-Q4" sheet and before "总计".
# 2. Insert a new top data-row into "总计" summarising the new quarter and
#    shift the existing rows' running index down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create "2022-Q1" sheet positioned right after "2021-Q4"
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $q4)
$newSheet.Name = "2022-Q1"

# Match page margins used by the other quarterly detail sheets
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Match the sheetPr/outlinePr defaults (summaryBelow/summaryRight) carried
# by all the sibling sheets.
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1

# Seed header row B1:H1 by copying formatting (bold / border / center-top)
# from the "2021-Q4" header, then overwrite the text.
$q4.Range("B1:H1").Copy($newSheet.Range("B1"))

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Fund holding rows. Column A keeps the running 0-based index (styled like
# the other sheets' index column); columns B-G are free-text (fund code /
# name / figures are reported as text strings, not numbers); column H is a
# plain integer rank.
$rows = @(
    @{ idx = 0; code = "910004"; name = "东方红启恒三年持有期混合型证券投资基金A"; scale = "118.66"; pos = "86.48"; pct = "4.00"; mv = "4.7464"; rank = 9 },
    @{ idx = 1; code = "011724"; name = "东方红启恒三年持有期混合型证券投资基金B"; scale = "110.24"; pos = "86.48"; pct = "4.00"; mv = "4.4096"; rank = 9 },
    @{ idx = 2; code = "010059"; name = "东方红鼎元3个月定期开放混合";             scale = "27.56";  pos = "89.58"; pct = "3.55"; mv = "0.9784"; rank = 10 },
    @{ idx = 3; code = "910009"; name = "东方红启程三年持有期混合型证券投资基金A"; scale = "8.06";   pos = "87.95"; pct = "4.01"; mv = "0.3232"; rank = 9 },
    @{ idx = 4; code = "001521"; name = "国寿安保成长优选股票";                   scale = "4.24";   pos = "87.83"; pct = "4.24"; mv = "0.1798"; rank = 6 },
    @{ idx = 5; code = "002861"; name = "工银瑞信智能制造股票";                   scale = "1.34";   pos = "93.93"; pct = "4.33"; mv = "0.0580"; rank = 10 },
    @{ idx = 6; code = "008082"; name = "国寿安保研究精选混合A";                  scale = "0.52";   pos = "91.60"; pct = "5.49"; mv = "0.0285"; rank = 5 },
    @{ idx = 7; code = "008083"; name = "国寿安保研究精选混合C";                  scale = "0.15";   pos = "91.60"; pct = "5.49"; mv = "0.0082"; rank = 5 }
)

$r = 2
foreach ($row in $rows) {
    # Column A: 0-based index, reuse the index-column styling (s="2") from
    # an existing sheet via copy.
    $q4.Range("A2").Copy($newSheet.Cells.Item($r, 1))
    $newSheet.Cells.Item($r, 1).Value = $row.idx

    $dataRange = $newSheet.Range("B$r`:G$r")
    $dataRange.NumberFormat = "@"
    $newSheet.Cells.Item($r, 2).Value = $row.code
    $newSheet.Cells.Item($r, 3).Value = $row.name
    $newSheet.Cells.Item($r, 4).Value = $row.scale
    $newSheet.Cells.Item($r, 5).Value = $row.pos
    $newSheet.Cells.Item($r, 6).Value = $row.pct
    $newSheet.Cells.Item($r, 7).Value = $row.mv
    $dataRange.Style = "Normal"

    $newSheet.Cells.Item($r, 8).Value = $row.rank

    $r = $r + 1
}

# ---------------------------------------------------------------------
# Step 2: insert the new quarter's summary row at the top of "总计"
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Shift existing data rows down, preserving their own formatting/values.
$total.Rows(2).Insert()
$total.Range("B2:D2").Style = "Normal"

# Re-number the (now shifted) running index column A for the pre-existing
# rows: old 0..4 (rows 3..7) becomes 1..5. (NOTE: `.Value` only supports
# writes in this host - use `.Value2` to read the current number back.)
for ($row = 7; $row -ge 3; $row--) {
    $old = $total.Cells.Item($row, 1).Value2
    $total.Cells.Item($row, 1).Value = $old + 1
}

# New top row: 2022-Q1 summary. Grab the index-column style (s="2") from
# the row just below (still holding its pre-shift value/format) before we
# overwrite it.
$total.Range("A3").Copy($total.Cells.Item(2, 1))
$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 8
$total.Cells.Item(2, 4).Value = 10.73

# Restore the originally active sheet/tab (unchanged by this edit).
$wb.Worksheets.Item("2020-Q4").Activate()
